$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Mapping SRU - CPSV-AP"
Write-Host $ws.Name
